# analise descritiva - acrescimo das metas
# Adds the "meta" / "meta_avg" / "meta_std" / "meta_min" / "meta_max" summary columns
# (inserted right after "taxa_sucesso", pushing "arrecadado_sucesso.."maior_ano" five
# columns to the right) and appends 5 more header cells at the end of the header row
# for the columns this pushes past the previous edge of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K', 'L', 'M', 'N', 'O', 'P', 'Q', 'R', 'S', 'T', 'U', 'V', 'W', 'X', 'Y', 'Z', 'AA')

# --- Row 1: header labels, written left-to-right in final column order ---
$headerText = @('modalidade', 'origem', 'total', 'total_sucesso', 'particip', 'taxa_sucesso', 'meta', 'meta_avg', 'meta_std', 'meta_min', 'meta_max', 'arrecadado_sucesso', 'arrecadado_avg', 'arrecadado_std', 'arrecadado_min', 'arrecadado_max', 'apoio_medio', 'apoio_std', 'apoio_min', 'apoio_max', 'contribuicoes', 'contribuicoes_med', 'contribuicoes_std', 'contribuicoes_min', 'contribuicoes_max', 'menor_ano', 'maior_ano')
for ($i = 0; $i -lt $headerText.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headerText[$i]
}
$hdrRange = $ws.Range("A1:AA1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160
$hdrRange.Borders.LineStyle = 1

# --- Row 2 ---
$row2Vals = @('flex', 'apoia.se', 5, 0, 0.003405994550408719, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row2Vals.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $row2Vals[$i]
}

# --- Row 3 ---
$row3Vals = @('flex', 'catarse', 1463, 1383, 0.9965940054495913, 0.9453178400546821, 15599716.7029188, 11279.62162177787, 16430.30708090436, 12.04441558726698, 198811.9434626772, 18362131.9375591, 13277.02960054888, 33934.82811955066, 10.77163914429046, 708972.7845446636, 77.41063997458096, 39.50983355883143, 10.77163914429046, 461.5197709071476, 203646, 147.2494577006508, 327.6748910926806, 1, 7954, 2016, 2023)
for ($i = 0; $i -lt $row3Vals.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3Vals[$i]
}

# --- Number formats for the data rows, matching the column groupings ---
$ws.Range("C2:D3").NumberFormat = "#,##0"
$ws.Range("E2:F3").NumberFormat = "0.00%"
$ws.Range("G2:T3").NumberFormat = "R$ #,##0.00"
$ws.Range("U2:Y3").NumberFormat = "#,##0"
$ws.Range("Z2:AA3").NumberFormat = "General"
